$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing contents so the shared-string table is rebuilt from
# scratch in the exact order the new values are written below.
$ws.Range("A1:AC10").ClearContents()

$data = @(
    @("Staff","Day 1","Day 2","Day 3","Day 4","Day 5","Day 6","Day 7","Day 8","Day 9","Day 10","Day 11","Day 12","Day 13","Day 14","Day 15","Day 16","Day 17","Day 18","Day 19","Day 20","Day 21","Day 22","Day 23","Day 24","Day 25","Day 26","Day 27","Day 28"),
    @("Staff_1","M1","M1","M3","M1","DO","M1","M1","A1","DO","M1","M1","M1","M1","M3","A1","A1","M3","M1","DO","M1","M1","A1","DO","PH","PH","M3","M1","M1"),
    @("Staff_2","M2","A2","M2","A2","M1","M1","DO","A2","M1","A1","M2","DO","M2","A2","M1","M2","DO","M2","A1","M2","A2","M2","M2","PH","PH","DO","M2","A2"),
    @("Staff_3","A1","DO","M3","M1","A1","A1","M1","DO","M1","M1","A1","A1","M1","M3","DO","M1","A1","M1","M3","A1","M1","DO","M1","PH","PH","A1","M1","M3"),
    @("Staff_4","DO","M2","A2","M1","M2","M1","A2","M1","A2","M2","DO","M2","M2","M1","M2","DO","M2","A2","M2","M1","A1","M2","M2","PH","PH","DO","M2","M2"),
    @("Staff_5","DO","M2","A2","M2","M1","M1","A2","M2","M2","DO","M2","M2","M1","A1","M2","DO","M2","A1","M1","M2","A2","M2","DO","PH","PH","M2","M2","A2"),
    @("Staff_6","A1","A1","DO","A1","A1","A1","M3","A1","A1","DO","A1","A1","A1","M3","A1","A1","A1","DO","A1","A1","M3","A1","A1","PH","PH","M3","A1","DO"),
    @("Staff_7","A2","DO","M2","A2","A2","M1","A1","M2","A2","A2","DO","M1","M1","A2","A2","A2","M2","A2","A1","M1","DO","A2","A2","PH","PH","DO","M2","A2"),
    @("Staff_8","M2","M2","A2","DO","M1","A2","M1","M2","M2","A2","A2","M1","A1","DO","M2","M2","DO","M2","M2","A1","A1","M2","M2","PH","PH","A2","A2","DO"),
    @("Staff_9","M2","A2","A2","M1","A2","M1","DO","A2","A2","M2","M2","A1","A1","DO","A2","M2","A2","A1","DO","M2","M1","A2","A2","PH","PH","M2","A2","DO")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

$wb.Save()
